$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header for new column F, copying the format from E1 (bold header style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "RXNO_DEF"

# Fill in RXNO_DEF definitions for each row
$ws.Range("F2").Value = '[''A generically dependent continuant that is about some thing. [IAO]'']'
$ws.Range("F3").Value = '[''p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]'', locstr("Process, i.e., a physical entity with a temporal evolution that ''has a meaning for the ontologist''", ''en'')]'
$ws.Range("F4").Value = '[''B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]'']'
$ws.Range("F5").Value = '[''Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]'']'
$ws.Range("F6").Value = '[''An elemental molecule consisting of two bivalently-bonded oxygen atoms. [Allotrope]'']'
$ws.Range("F7").Value = '[''An elemental molecule consisting of two trivalently-bonded nitrogen atoms. [CHEBI]'']'
$ws.Range("F8").Value = '[''A one-carbon compound with formula CO2 in which the carbon is attached to each oxygen atom by a double bond. [CHEBI]'']'
$ws.Range("F9").Value = '[''An azane that consists of a single nitrogen atom covelently bonded to three hydrogen atoms. [CHEBI]'']'
$ws.Range("F10").Value = '[''A one-carbon compound in which the carbon is joined only to a single oxygen. [CHEBI]'']'
$ws.Range("F11").Value = '[''A processual entity that realizes a plan which is the concretization of a plan specification. [IAO]'']'
$ws.Range("F12").Value = '[''Organonitrogen compounds that are derivatives of isocyanic acid; compounds containing the isocyanate functional group ‒N=C=O (as opposed to the cyanate group, -O-C≡N). [CHEBI]'']'
$ws.Range("F13").Value = '[''A directive information entity that describes an intended process endpoint. When part of a plan specification the concretization is realized in a planned process in which the bearer tries to effect the world so that the process endpoint is achieved. [IAO]'']'
$ws.Range("F14").Value = '[]'
$ws.Range("F15").Value = '[]'

$excel.CutCopyMode = 0
